$d = $word.ActiveDocument

# Useful special characters
$lq   = [char]0x201C   # left double quotation mark  "
$rq   = [char]0x201D   # right double quotation mark "
$ndash = [char]0x2013  # en dash –
$apos = [char]0x2019   # right single quotation mark '

# --------------------------------------------------------------------
# 1. Achievement heading:
#    "Achievement in Delivering SaaS, Data-Intensive, Automation and
#     Analytics Applications"
#    -> "Achievement in Delivering Application Modernization, Data
#        Integration, Automation and Analytics Initiatives"
#    (Find text starts at "Delivering" rather than "Achievement in" so
#     the merged run inherits the plain formatting of the "De…" run
#     instead of the "Achievement in " run's extra eastAsia/lang attrs.)
# --------------------------------------------------------------------
$d.Content.Find.Execute(
    "Delivering SaaS, Data-Intensive, Automation and Analytics Applications",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Delivering Application Modernization, Data Integration, Automation and Analytics Initiatives",
    2) | Out-Null

# --------------------------------------------------------------------
# 2 & 3. Bullet about the modernized platform:
#    "Consolidated 7 risk data areas for 5 asset classes into 20 key
#     reports from internal datawarehouse and Sungard risk platform"
#    -> "Modernized application platform for 7 risk areas over 5 asset
#        classes from internal datawarehouse to Sungard risk system"
# --------------------------------------------------------------------
$d.Content.Find.Execute(
    "Consolidated 7 risk data areas for 5 asset classes into 20 key  reports from internal datawarehouse and Sungard risk platform",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Modernized application platform for 7 risk areas over 5 asset classes from internal datawarehouse to Sungard risk system",
    2) | Out-Null

# --------------------------------------------------------------------
# 4. "Delivered bank's internal repository of 200 SOX controls with
#     assessments and testing evidences from 500 key processes"
#    -> "Delivered bankwide self-assessment with repository of 200 SOX
#        controls, assessments and testing from 500 key processes"
# --------------------------------------------------------------------
$d.Content.Find.Execute(
    ("Delivered bank" + $apos + "s internal repository of 200 SOX controls with assessments and testing evidences from 500 key processes"),
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Delivered bankwide self-assessment with repository of 200 SOX controls, assessments and testing from 500 key processes",
    2) | Out-Null

# --------------------------------------------------------------------
# 5. Merge the split "(methodology and audit findings fully disclosed
#    to external auditors)" runs into one (no visible text change).
# --------------------------------------------------------------------
$d.Content.Find.Execute(
    "(methodology and audit findings fully disclosed to external auditors)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(methodology and audit findings fully disclosed to external auditors)",
    2) | Out-Null

# --------------------------------------------------------------------
# 6. Drop the old "_GoBack" bookmark that sits right after "/PEGA"
#    (a no-op text replace whose range spans the bookmark removes it,
#    same as typing over a bookmarked position in real Word).
# --------------------------------------------------------------------
$d.Content.Find.Execute(
    "BLUEPRISM/PEGA RPA,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "BLUEPRISM/PEGA RPA,",
    2) | Out-Null

# --------------------------------------------------------------------
# 7. "RPA, APPIAN Treasury Workflow, CI/CD  AZURE DEVOPS"
#    -> "RPA, APPIAN Treasury Workflow, CI/CD AZURE TFS, SONIC,
#         SAILPOINT, SPLUNK"
# --------------------------------------------------------------------
$d.Content.Find.Execute(
    "RPA, APPIAN Treasury Workflow, CI/CD  AZURE DEVOPS",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "RPA, APPIAN Treasury Workflow, CI/CD AZURE TFS, SONIC, SAILPOINT, SPLUNK",
    2) | Out-Null

# --------------------------------------------------------------------
# 8. Rewrite the "Tactical Solution" bullet and move the "_GoBack"
#    bookmark to right after the newly-added "...bank's final " text.
# --------------------------------------------------------------------
$findText8 = "Developed the " + $lq + "Internal Control Repository Tactical Solution" + $rq + " " + $ndash + " a SQL Server-based application to store the assessment (test, deficiencies tracking, remediation) and assertion (exceptions, sign-off) of SOX controls from 200 business unit"
$replaceText8 = "Developed the Tactical Solution to store the assessment (test, deficiencies tracking, remediation) and assertion (exceptions, sign-off) of SOX controls from 200 business unit. This was the prototype for the bank" + $apos + "s final CSA & Internal Control Repository"

$d.Content.Find.Execute($findText8, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText8, 2) | Out-Null

$r8 = $d.Content
$r8.Find.Execute(("bank" + $apos + "s final "), $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmRange = $d.Range($r8.End, $r8.End)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
